# Adds 5 new days of mobility data (columns MD:MH) to rows 1-5, extending
# the existing date series (row 1) and the per-country numeric series
# (rows 2-5), matching the author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds dates formatted with the workbook's existing date style.
# Copy the format from the last populated date cell (MC1) onto the new
# cells first, so the new cells reuse the same style index instead of
# Excel minting a duplicate cellXf, then fill in the new date values.
$ws.Range("MC1").Copy()
$ws.Range("MD1:MH1").PasteSpecial(-4122)

$ws.Range("MD1").Value = 44178
$ws.Range("ME1").Value = 44179
$ws.Range("MF1").Value = 44180
$ws.Range("MG1").Value = 44181
$ws.Range("MH1").Value = 44182

# Row 2 new values
$ws.Range("MD2").Value = 43.11
$ws.Range("ME2").Value = 58.29
$ws.Range("MF2").Value = 55.5
$ws.Range("MG2").Value = 53.85
$ws.Range("MH2").Value = 56.71

# Row 3 new values
$ws.Range("MD3").Value = 46.24
$ws.Range("ME3").Value = 42.21
$ws.Range("MF3").Value = 45.59
$ws.Range("MG3").Value = 36.99
$ws.Range("MH3").Value = 37.18

# Row 4 new values
$ws.Range("MD4").Value = 52.07
$ws.Range("ME4").Value = 60.35
$ws.Range("MF4").Value = 61.9
$ws.Range("MG4").Value = 60.74
$ws.Range("MH4").Value = 59.07

# Row 5 new values
$ws.Range("MD5").Value = 23.18
$ws.Range("ME5").Value = 32.19
$ws.Range("MF5").Value = 29.06
$ws.Range("MG5").Value = 30.1
$ws.Range("MH5").Value = 29.66
